$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28 (pushes the existing rows 28-36 down to 29-37)
$ws.Rows("28:28").Insert()

# Populate the new row with the added radio source entry
$ws.Range("A28").Value = "J1530+1049"
$ws.Range("B28").Value = 5.72
$ws.Range("C28").Value = 174.88
$ws.Range("D28").Value = -1.4
$ws.Range("E28").Value = "Saxena+18"

# Extend the sorted range/condition to include the new row, keeping data order intact
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B1:B34"), 0, 2)
$ws.Sort.SetRange($ws.Range("A3:G34"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Match the saved view state: scrolled/selected cell after the edit
$ws.Range("E29").Select()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
